# "updated thread local driver" - the LoginStatus.xlsx test-data sheet gets
# two more worker-thread status columns (K, L) appended, mirroring the
# existing C..J "status/PASS/FAIL" columns exactly (same header text, same
# per-row values, same grey header fill, same column sizing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source column J holds the canonical pattern we need to replicate onto K, L.
$sourceCol = 10   # column J
$newCols = @(11, 12)  # columns K, L

foreach ($col in $newCols) {
    $srcCell = $ws.Cells.Item(1, $sourceCol)
    $dstHeader = $ws.Cells.Item(1, $col)

    # Header cell: same text ("status") and same grey (indexed 55) fill used
    # by every other header cell in the row.
    $dstHeader.Value = $srcCell.Value()
    $dstHeader.Interior.ColorIndex = 48

    # Data rows 2-6: replicate PASS/PASS/FAIL/FAIL/FAIL (no explicit style,
    # same as column J).
    for ($r = 2; $r -le 6; $r++) {
        $src = $ws.Cells.Item($r, $sourceCol)
        $dst = $ws.Cells.Item($r, $col)
        $dst.Value = $src.Value()
    }

    # Match column width/sizing to the source column.
    $ws.Columns.Item($col).ColumnWidth = 5.43
}
